# Apply the "Velmi dobra shoda reakci" update:
#  - add a new row of reaction-comparison results (row 5)
#  - apply the "Carka" (thousands/comma) number style to the newly
#    compared F/G/H (and I/J where relevant) cells in rows 1-3
#  - move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

$commaFormat = '_-* #,##0.00\ _K_č_-;\-* #,##0.00\ _K_č_-;_-* "-"??\ _K_č_-;_-@_-'

# Rows 1-3: give the previously unformatted F/G/H (and I3/J3) cells the
# same "Carka" number format already used by columns B:E so they read s="1".
$ws.Range("F1:I1").NumberFormat = $commaFormat
$ws.Range("F2:H2").NumberFormat = $commaFormat
$ws.Range("F3:J3").NumberFormat = $commaFormat

# Row 5: new set of results, matching the layout of rows 1-3.
$ws.Range("A5").Value = "PRJ-6076_DLC24_OVS_V11.0_N_S1_FATIGUE"
$ws.Range("B5").Value = 2104401.16788114
$ws.Range("C5").Value = 1233079.4681556099
$ws.Range("D5").Value = 1401312.69478258
$ws.Range("E5").Value = 748528.80624845705
$ws.Range("F5").Value = 386
$ws.Range("G5").Value = 267
$ws.Range("H5").Value = 4201
$ws.Range("I5").Value = 140
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 7

# Move the active cell/selection as it was left in the saved workbook.
$ws.Range("H12").Select() | Out-Null
